# Re-order the currency-file summary rows (A:D); the E (built_in_total)
# column is left untouched per-row, matching the source change.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "87811004_1121_MX"
$ws.Cells.Item(2, 2).Value = 78
$ws.Cells.Item(2, 3).Value = "MXN"
$ws.Cells.Item(2, 4).Value = 7567.7

$ws.Cells.Item(3, 1).Value = "87811004_1121_BR"
$ws.Cells.Item(3, 2).Value = 38
$ws.Cells.Item(3, 3).Value = "BRL"
$ws.Cells.Item(3, 4).Value = 440.44

$ws.Cells.Item(4, 1).Value = "87811004_1121_CA"
$ws.Cells.Item(4, 2).Value = 287
$ws.Cells.Item(4, 3).Value = "CAD"
$ws.Cells.Item(4, 4).Value = 1332.8

$ws.Cells.Item(5, 1).Value = "87811004_1121_LL"
$ws.Cells.Item(5, 2).Value = 43
$ws.Cells.Item(5, 3).Value = "USD"
$ws.Cells.Item(5, 4).Value = 112

$ws.Cells.Item(6, 1).Value = "87811004_1121_BG"
$ws.Cells.Item(6, 2).Value = 6
$ws.Cells.Item(6, 3).Value = "BGN"
$ws.Cells.Item(6, 4).Value = 16.04

$ws.Cells.Item(7, 1).Value = "87811004_1121_HU"
$ws.Cells.Item(7, 2).Value = 900
$ws.Cells.Item(7, 3).Value = "HUF"
$ws.Cells.Item(7, 4).Value = 1915878

$ws.Cells.Item(8, 1).Value = "87811004_1121_PE"
$ws.Cells.Item(8, 2).Value = 21
$ws.Cells.Item(8, 3).Value = "PEN"
$ws.Cells.Item(8, 4).Value = 174.3

$ws.Cells.Item(9, 1).Value = "87811004_1121_RO"
$ws.Cells.Item(9, 2).Value = 494
$ws.Cells.Item(9, 3).Value = "RON"
$ws.Cells.Item(9, 4).Value = 9795.57

$ws.Cells.Item(10, 1).Value = "87811004_1121_EU"
$ws.Cells.Item(10, 2).Value = 746
$ws.Cells.Item(10, 3).Value = "EUR"
$ws.Cells.Item(10, 4).Value = 3309.31

$ws.Cells.Item(11, 1).Value = "87811004_1121_PL"
$ws.Cells.Item(11, 2).Value = 42
$ws.Cells.Item(11, 3).Value = "PLN"
$ws.Cells.Item(11, 4).Value = 401.94

$ws.Cells.Item(12, 1).Value = "87811004_1121_DK"
$ws.Cells.Item(12, 2).Value = 22
$ws.Cells.Item(12, 3).Value = "DKK"
$ws.Cells.Item(12, 4).Value = 440.16

$ws.Cells.Item(13, 1).Value = "87811004_1121_GB"
$ws.Cells.Item(13, 2).Value = 400
$ws.Cells.Item(13, 3).Value = "GBP"
$ws.Cells.Item(13, 4).Value = 1309.43

$ws.Cells.Item(14, 1).Value = "87811004_1121_SE"
$ws.Cells.Item(14, 2).Value = 36
$ws.Cells.Item(14, 3).Value = "SEK"
$ws.Cells.Item(14, 4).Value = 919.87

$ws.Cells.Item(15, 1).Value = "87811004_1121_JP"
$ws.Cells.Item(15, 2).Value = 23
$ws.Cells.Item(15, 3).Value = "JPY"
$ws.Cells.Item(15, 4).Value = 7546

$ws.Cells.Item(16, 1).Value = "87811004_1121_CO"
$ws.Cells.Item(16, 2).Value = 31
$ws.Cells.Item(16, 3).Value = "COP"
$ws.Cells.Item(16, 4).Value = 347830

$ws.Cells.Item(17, 1).Value = "87811004_1121_CL"
$ws.Cells.Item(17, 2).Value = 35
$ws.Cells.Item(17, 3).Value = "CLP"
$ws.Cells.Item(17, 4).Value = 73248

$ws.Cells.Item(18, 1).Value = "87811004_1121_CZ"
$ws.Cells.Item(18, 2).Value = 23
$ws.Cells.Item(18, 3).Value = "CZK"
$ws.Cells.Item(18, 4).Value = 1718.21

$ws.Cells.Item(19, 1).Value = "87811004_1121_US"
$ws.Cells.Item(19, 2).Value = 1415
$ws.Cells.Item(19, 3).Value = "USD"
$ws.Cells.Item(19, 4).Value = 8030.4

$ws.Cells.Item(20, 1).Value = "87811004_1121_CH"
$ws.Cells.Item(20, 2).Value = 78
$ws.Cells.Item(20, 3).Value = "CHF"
$ws.Cells.Item(20, 4).Value = 387.56

$ws.Cells.Item(21, 1).Value = "87811004_1121_NO"
$ws.Cells.Item(21, 2).Value = 26
$ws.Cells.Item(21, 3).Value = "NOK"
$ws.Cells.Item(21, 4).Value = 655.2

$ws.Cells.Item(22, 1).Value = "87811004_1121_NZ"
$ws.Cells.Item(22, 2).Value = 34
$ws.Cells.Item(22, 3).Value = "NZD"
$ws.Cells.Item(22, 4).Value = 147.6

$ws.Cells.Item(23, 1).Value = "87811004_1121_AU"
$ws.Cells.Item(23, 2).Value = 307
$ws.Cells.Item(23, 3).Value = "AUD"
$ws.Cells.Item(23, 4).Value = 1617.46

